$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "REFS" -> "RELATIONS"
$ws.Range("B1").Value = "RELATIONS"

# Update "Parent: REQ-001" text (remove trailing newline)
$ws.Range("B3").Value = "Parent: REQ-001"

# Update "Parent: REQ-002" -> multi-line combined parent text
$ws.Range("B4").Value = "Parent: REQ-001" + [char]10 + "----" + [char]10 + "Parent: REQ-002"
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# Widen column B
$ws.Columns.Item(2).ColumnWidth = 36.7109375

# Update the table column name to match the new header
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item(2).Name = "RELATIONS"
